$d = $word.ActiveDocument

# 1. Replace the title runs ("Standard Work Procedure" + "Change 1") with a
#    single run of text "Change to title<<MARK>>" - the trailing marker
#    lets us anchor the "_GoBack" bookmark exactly at the end of the real
#    text (immediately after the run) before trimming the marker back off,
#    leaving a clean, unsplit run that matches "Change to title".
$marker = "<<MARK>>"
$d.Content.Find.Execute("Standard Work ProcedureChange 1", $true, $false, $false, $false, $false, $true, 1, $false, "Change to title" + $marker, 2)

$titlePara = $d.Paragraphs(1).Range
$titleText = $titlePara.Text
$markerLen = $marker.Length
$boundary = $titlePara.Start + $titleText.Length - $markerLen - 1

# 2. Move the "_GoBack" bookmark (previously sitting in a later, now-unused
#    paragraph) to just after the title text, inside paragraph 1. Adding a
#    bookmark with a name that already exists relocates it instead of
#    creating a duplicate, so the old bookmarkStart/bookmarkEnd pair goes
#    away on its own.
$bmRange = $d.Range($boundary, $boundary)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3. Strip the temporary marker text back out, leaving "Change to title".
$markerRange = $d.Range($boundary, $boundary + $markerLen)
$markerRange.Delete()
